$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add new column I with header "hello" and a few numeric values in rows 2-4,
# matching the "Species removed when training and predicting" feature column.
$ws.Range("I1").Value = "hello"
$ws.Range("I2").Value = 1
$ws.Range("I3").Value = 1
$ws.Range("I4").Value = 1
